$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Scene events that trigger a battle (monster-fight rows) now support
#     being triggered multiple times: TriggerMulti column (C) 0 -> 1.
$ws.Range("C5:C10").Value = 1

# --- New row 24: "portal" scene event (传送门 / portal) ---------------------
$ws.Cells.Item(24, 1).Value = 42010015
$ws.Cells.Item(24, 2).Value = "传送门"
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = "portal"
$ws.Cells.Item(24, 5).Value = "portal"
$ws.Cells.Item(24, 6).Value = "portal"

# --- Punish/reward events now allow their "identifier" text to be reused
#     under a distinct entry (sandpile / mushroom) -----------------------
$ws.Cells.Item(18, 6).Value = "sandpile"
$ws.Cells.Item(19, 6).Value = "mushroom"

# --- Treasure box event's Ename is renamed from "treasurebox" to "treasure"
$ws.Cells.Item(21, 4).Value = "treasure"

# --- New row 25: "falling stone" scene event (落石 / stone) ----------------
$ws.Cells.Item(25, 1).Value = 42010016
$ws.Cells.Item(25, 2).Value = "落石"
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = "stone"
$ws.Cells.Item(25, 5).Value = "stone"
$ws.Cells.Item(25, 6).Value = "stone"
$ws.Cells.Item(25, 19).Value = 80

# --- Grow the "表3" table/autofilter to cover the two new rows -------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:X25"))

# --- Match the saved selection/active cell from the source workbook --------
[void]$ws.Range("S25").Select()
